$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B12").Value = "positive regulation of potassium ion transmembrane transporter activity"
$ws.Range("B13").Value = "positive regulation of cation channel activity"
$ws.Range("B14").Value = "ventricular cardiac muscle cell action potential"
$ws.Range("B15").Value = "membrane depolarization during SA node cell action potential"
$ws.Range("B16").Value = "regulation of atrial cardiac muscle cell action potential"
$ws.Range("B17").Value = "regulation of heart rate by cardiac conduction"
$ws.Range("B18").Value = "regulation of SA node cell action potential"

$ws.Range("B26").Value = "regulation of membrane permeability"
$ws.Range("B27").Value = "negative regulation of vernalization response"
$ws.Range("B28").Value = "response to herbivore"
$ws.Range("B29").Value = "response to molecule of fungal origin"
$ws.Range("B30").Value = "protein localization to cell surface"

$ws.Range("B81").Value = "retinal rod cell development"
$ws.Range("B82").Value = "photoreceptor cell outer segment organization"
